$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 138.02942
$ws.Range("I9").Value = 114.129036
$ws.Range("K9").Value = 114.129036
$ws.Range("M9").Value = 54.870964
$ws.Range("H17").Value = 759.34045
$ws.Range("J17").Value = 766.2826
$ws.Range("L17").Value = 2298.8478
$ws.Range("N17").Value = -2634.8478
$ws.Range("H57").Value = 63699
$ws.Range("J57").Value = 63699
$ws.Range("L57").Value = 191097
$ws.Range("N57").Value = -192095
$ws.Range("H70").Value = 9314.429
$ws.Range("I70").Value = 3714.2856
$ws.Range("J70").Value = 14914.571
$ws.Range("K70").Value = 11142.8568
$ws.Range("L70").Value = 44743.713
$ws.Range("M70").Value = -10872.8568
$ws.Range("N70").Value = -45283.713
$ws.Range("H73").Value = 9314.429
$ws.Range("I73").Value = 3714.2856
$ws.Range("J73").Value = 14914.571
$ws.Range("K73").Value = 11142.8568
$ws.Range("L73").Value = 44743.713
$ws.Range("M73").Value = -10206.8568
$ws.Range("N73").Value = -46615.713
$ws.Range("H86").Value = 4588.2
$ws.Range("I86").Value = 4419
$ws.Range("K86").Value = 4419
$ws.Range("M86").Value = -3296
$ws.Range("H89").Value = 4588.2
$ws.Range("I89").Value = 4419
$ws.Range("K89").Value = 22095
$ws.Range("M89").Value = -16479
$ws.Range("H98").Value = 1215.3489
$ws.Range("I98").Value = 1095.359
$ws.Range("J98").Value = 2385.25
$ws.Range("K98").Value = 1095.359
$ws.Range("L98").Value = 2385.25
$ws.Range("M98").Value = 402.6410000000001
$ws.Range("N98").Value = -5381.25
$ws.Range("I113").Value = 2830.5
$ws.Range("J113").Value = 9802.6
$ws.Range("K113").Value = 2830.5
$ws.Range("L113").Value = 9802.6
$ws.Range("M113").Value = 423.5
$ws.Range("N113").Value = -16310.6
$ws.Range("H116").Value = 2733.1667
$ws.Range("I116").Value = 2700
$ws.Range("J116").Value = 2799.5
$ws.Range("K116").Value = 2700
$ws.Range("L116").Value = 2799.5
$ws.Range("M116").Value = 742
$ws.Range("N116").Value = -9683.5
$ws.Range("H122").Value = 1215.3489
$ws.Range("I122").Value = 1095.359
$ws.Range("J122").Value = 2385.25
$ws.Range("K122").Value = 3286.077
$ws.Range("L122").Value = 7155.75
$ws.Range("M122").Value = -836.0769999999998
$ws.Range("N122").Value = -12055.75
$ws.Range("H131").Value = 1164.9131
$ws.Range("I131").Value = 744.7368
$ws.Range("J131").Value = 3160.75
$ws.Range("K131").Value = 2234.2104
$ws.Range("L131").Value = 9482.25
$ws.Range("M131").Value = 2805.7896
$ws.Range("N131").Value = -19562.25
$ws.Range("H137").Value = 2598.7585
$ws.Range("I137").Value = 2643.111
$ws.Range("K137").Value = 7929.333
$ws.Range("M137").Value = -5379.333
$ws.Range("H138").Value = 2061562.1
$ws.Range("I138").Value = 989.5714
$ws.Range("K138").Value = 2968.7142
$ws.Range("M138").Value = 2171.2858

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H32").Value = 1458.8549
$ws.Range("I32").Value = 1208.098
$ws.Range("J32").Value = 2621.4546
$ws.Range("K32").Value = 1208.098
$ws.Range("L32").Value = 2621.4546
$ws.Range("M32").Value = -921.098
$ws.Range("N32").Value = -3195.4546
$ws.Range("H45").Value = 2072.4
$ws.Range("J45").Value = 2999
$ws.Range("L45").Value = 2999
$ws.Range("N45").Value = -3753
$ws.Range("H61").Value = 9365.134
$ws.Range("I61").Value = 6043.5454
$ws.Range("K61").Value = 6043.5454
$ws.Range("M61").Value = -5831.5454
$ws.Range("H70").Value = 85000
$ws.Range("J70").Value = 85000
$ws.Range("L70").Value = 85000
$ws.Range("N70").Value = -85540
$ws.Range("H73").Value = 85000
$ws.Range("J73").Value = 85000
$ws.Range("L73").Value = 85000
$ws.Range("N73").Value = -86872
$ws.Range("H102").Value = 8001452.5
$ws.Range("I102").Value = 1273.2667
$ws.Range("K102").Value = 1273.2667
$ws.Range("M102").Value = 348.7333000000001
$ws.Range("H122").Value = 2510.6667
$ws.Range("I122").Value = 2513.4443
$ws.Range("J122").Value = 2506.5
$ws.Range("K122").Value = 7540.3329
$ws.Range("L122").Value = 7519.5
$ws.Range("M122").Value = -5090.3329
$ws.Range("N122").Value = -12419.5
$ws.Range("H132").Value = 5235.857
$ws.Range("I132").Value = 4274
$ws.Range("J132").Value = 8121.4287
$ws.Range("K132").Value = 12822
$ws.Range("L132").Value = 24364.2861
$ws.Range("M132").Value = -10292
$ws.Range("N132").Value = -29424.2861
$ws.Range("H136").Value = 9365.134
$ws.Range("I136").Value = 6043.5454
$ws.Range("K136").Value = 18130.6362
$ws.Range("M136").Value = -15580.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 36750
$ws.Range("J55").Value = 22500
$ws.Range("L55").Value = 22500
$ws.Range("N55").Value = -23046
$ws.Range("H56").Value = 9554.777
$ws.Range("J56").Value = 9554.777
$ws.Range("L56").Value = 9554.777
$ws.Range("N56").Value = -11032.777
$ws.Range("H82").Value = 5000
$ws.Range("I82").Value = 5000
$ws.Range("K82").Value = 5000
$ws.Range("M82").Value = -4617
$ws.Range("H85").Value = 5000
$ws.Range("I85").Value = 5000
$ws.Range("K85").Value = 5000
$ws.Range("M85").Value = -3674
$ws.Range("H94").Value = 2388.5
$ws.Range("I94").Value = 2238.2173
$ws.Range("J94").Value = 3079.8
$ws.Range("K94").Value = 2238.2173
$ws.Range("L94").Value = 3079.8
$ws.Range("M94").Value = -1787.2173
$ws.Range("N94").Value = -3981.8
$ws.Range("H97").Value = 32436
$ws.Range("I97").Value = 20428
$ws.Range("J97").Value = 44444
$ws.Range("K97").Value = 20428
$ws.Range("L97").Value = 44444
$ws.Range("M97").Value = -19437
$ws.Range("N97").Value = -46426
$ws.Range("H99").Value = 2028.1428
$ws.Range("I99").Value = 1430.625
$ws.Range("J99").Value = 3940.2
$ws.Range("K99").Value = 1430.625
$ws.Range("L99").Value = 3940.2
$ws.Range("M99").Value = 67.375
$ws.Range("N99").Value = -6936.2
$ws.Range("H105").Value = 4937.524
$ws.Range("I105").Value = 3656.2
$ws.Range("J105").Value = 6102.364
$ws.Range("K105").Value = 3656.2
$ws.Range("L105").Value = 6102.364
$ws.Range("M105").Value = -1909.2
$ws.Range("N105").Value = -9596.364
$ws.Range("H107").Value = 2078.4167
$ws.Range("I107").Value = 2255.7144
$ws.Range("J107").Value = 837.3333
$ws.Range("K107").Value = 2255.7144
$ws.Range("L107").Value = 837.3333
$ws.Range("M107").Value = -335.7143999999998
$ws.Range("N107").Value = -4677.3333
$ws.Range("H122").Value = 87000
$ws.Range("J122").Value = 87000
$ws.Range("L122").Value = 87000
$ws.Range("N122").Value = -96800
$ws.Range("H134").Value = 5341
$ws.Range("I134").Value = 5450.769
$ws.Range("K134").Value = 16352.307
$ws.Range("M134").Value = -13817.307

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3384.75
$ws.Range("I31").Value = 2704.7334
$ws.Range("J31").Value = 3984.7646
$ws.Range("K31").Value = 2704.7334
$ws.Range("L31").Value = 3984.7646
$ws.Range("M31").Value = -2409.7334
$ws.Range("N31").Value = -4574.7646
$ws.Range("H34").Value = 3384.75
$ws.Range("I34").Value = 2704.7334
$ws.Range("J34").Value = 3984.7646
$ws.Range("K34").Value = 2704.7334
$ws.Range("L34").Value = 3984.7646
$ws.Range("M34").Value = -2502.7334
$ws.Range("N34").Value = -4388.7646
$ws.Range("H99").Value = 4586.8076
$ws.Range("J99").Value = 6387.5
$ws.Range("L99").Value = 6387.5
$ws.Range("N99").Value = -9383.5
$ws.Range("H122").Value = 2115.913
$ws.Range("I122").Value = 1395.5
$ws.Range("J122").Value = 4709.4
$ws.Range("K122").Value = 4186.5
$ws.Range("L122").Value = 14128.2
$ws.Range("M122").Value = -1736.5
$ws.Range("N122").Value = -19028.2
$ws.Range("H126").Value = 4586.8076
$ws.Range("J126").Value = 6387.5
$ws.Range("L126").Value = 19162.5
$ws.Range("N126").Value = -24102.5
$ws.Range("H132").Value = 1999.8
$ws.Range("I132").Value = 1899.6666
$ws.Range("K132").Value = 5698.9998
$ws.Range("M132").Value = -3168.9998
$ws.Range("H134").Value = 1743.7333
$ws.Range("I134").Value = 1554.75
$ws.Range("K134").Value = 4664.25
$ws.Range("M134").Value = -2129.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H32").Value = 1000000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3000000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3000566
$ws.Range("H39").Value = 4234.4
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 45000
$ws.Range("N39").Value = -45588
$ws.Range("H55").Value = 8193.125
$ws.Range("J55").Value = 12680
$ws.Range("L55").Value = 38040
$ws.Range("N55").Value = -38394
$ws.Range("H121").Value = 4256.8184
$ws.Range("J121").Value = 4347.857
$ws.Range("L121").Value = 13043.571
$ws.Range("N121").Value = -15663.571
$ws.Range("H126").Value = 2861.25
$ws.Range("I126").Value = 2861.25
$ws.Range("K126").Value = 8583.75
$ws.Range("M126").Value = -3643.75
$ws.Range("H131").Value = 1584.4546
$ws.Range("J131").Value = 2006.5834
$ws.Range("L131").Value = 6019.7502
$ws.Range("N131").Value = -16099.7502
$ws.Range("H134").Value = 2667.3635
$ws.Range("I134").Value = 2667.3635
$ws.Range("K134").Value = 8002.0905
$ws.Range("M134").Value = -2932.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 1873.5
$ws.Range("I46").Value = 1873.5
$ws.Range("K46").Value = 1873.5
$ws.Range("M46").Value = -1717.5
$ws.Range("H92").Value = 24400
$ws.Range("J92").Value = 24400
$ws.Range("L92").Value = 24400
$ws.Range("N92").Value = -28144
$ws.Range("H122").Value = 62503332
$ws.Range("I122").Value = 71432024
$ws.Range("K122").Value = 214296072
$ws.Range("M122").Value = -214293622
$ws.Range("H126").Value = 33928.145
$ws.Range("I126").Value = 4582.8335
$ws.Range("K126").Value = 13748.5005
$ws.Range("M126").Value = -11278.5005
$ws.Range("H132").Value = 3791.1667
$ws.Range("I132").Value = 4149.625
$ws.Range("J132").Value = 1998.875
$ws.Range("K132").Value = 12448.875
$ws.Range("L132").Value = 5996.625
$ws.Range("M132").Value = -9918.875
$ws.Range("N132").Value = -11056.625
$ws.Range("H136").Value = 73750
$ws.Range("J136").Value = 73750
$ws.Range("L136").Value = 221250
$ws.Range("N136").Value = -226350

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5049.375
$ws.Range("I9").Value = 4342.143
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 4342.143
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = -4118.143
$ws.Range("N9").Value = -10448
$ws.Range("H40").Value = 6612.6665
$ws.Range("I40").Value = 6553.115
$ws.Range("J40").Value = 6999.75
$ws.Range("K40").Value = 6553.115
$ws.Range("L40").Value = 6999.75
$ws.Range("M40").Value = -6417.115
$ws.Range("N40").Value = -7271.75
$ws.Range("H46").Value = 4615.304
$ws.Range("J46").Value = 10511
$ws.Range("L46").Value = 10511
$ws.Range("N46").Value = -10887
$ws.Range("H69").Value = 88931.5
$ws.Range("J69").Value = 88931.5
$ws.Range("L69").Value = 88931.5
$ws.Range("N69").Value = -90553.5
$ws.Range("H72").Value = 88931.5
$ws.Range("J72").Value = 88931.5
$ws.Range("L72").Value = 266794.5
$ws.Range("N72").Value = -274906.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 4190
$ws.Range("I100").Value = 3980
$ws.Range("K100").Value = 3980
$ws.Range("M100").Value = -3439
$ws.Range("H122").Value = 4685.852
$ws.Range("I122").Value = 4012.75
$ws.Range("J122").Value = 5664.909
$ws.Range("K122").Value = 12038.25
$ws.Range("L122").Value = 16994.727
$ws.Range("M122").Value = -9588.25
$ws.Range("N122").Value = -21894.727
$ws.Range("H132").Value = 9666.333000000001
$ws.Range("I132").Value = 9000
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 27000
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -24470
$ws.Range("N132").Value = -35058.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 150056.42
$ws.Range("J81").Value = 10999
$ws.Range("L81").Value = 21998
$ws.Range("N81").Value = -24120
$ws.Range("H84").Value = 150056.42
$ws.Range("J84").Value = 10999
$ws.Range("L84").Value = 109990
$ws.Range("N84").Value = -120598
$ws.Range("H112").Value = 79989.5
$ws.Range("J112").Value = 79989.5
$ws.Range("L112").Value = 79989.5
$ws.Range("N112").Value = -82943.5
$ws.Range("H122").Value = 3291.524
$ws.Range("I122").Value = 3961.625
$ws.Range("J122").Value = 1147.2
$ws.Range("K122").Value = 11884.875
$ws.Range("L122").Value = 3441.6
$ws.Range("M122").Value = -9434.875
$ws.Range("N122").Value = -8341.6
$ws.Range("H132").Value = 1152
$ws.Range("I132").Value = 1004
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 3012
$ws.Range("L132").Value = 3900
$ws.Range("M132").Value = -482
$ws.Range("N132").Value = -8960
$ws.Range("H135").Value = 108000
$ws.Range("J135").Value = 108000
$ws.Range("L135").Value = 108000
$ws.Range("N135").Value = -118140
